$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# --- Settings sheet: new PythonLibraryPath / PythonPath / WorkingFolder rows ---
$ws1.Range("A7").Value = "PythonLibraryPath"
$ws1.Range("B7").Value = "C:\Users\Osmar\AppData\Local\Programs\Python\Python310\python310.dll"
$ws1.Range("A8").Value = "PythonPath"
$ws1.Range("B8").Value = "C:\Users\Osmar\AppData\Local\Programs\Python\Python310"
$ws1.Range("A9").Value = "WorkingFolder"
$ws1.Range("B9").Value = "D:\RoitRPA\RoitRPA"

# --- Constants sheet: MaxConsecutiveSystemExceptions 0 -> 2 ---
$ws2.Range("B3").Value = 2

# --- Constants sheet: ShouldMarkJobAsFaulted boolean 0 -> text "FALSE" ---
# A direct Value="FALSE" assignment gets auto-coerced to a Boolean by Excel,
# so round-trip it through a formula -> copy -> paste-values to land a real
# text shared-string in the cell (matching t="s" in the target file).
$ws2.Range("B17").Formula = '="FALSE"'
$ws2.Range("B17").Copy()
$ws2.Range("B17").PasteSpecial(-4163)

# --- Selection / active sheet bookkeeping ---
[void]$ws1.Range("B15").Select()
[void]$ws2.Activate()
[void]$ws2.Range("B18").Select()
